# Trade #89 closed at 2026-02-17 15:53:56 - unknown UNKNOWN +0.000%
#
# Updates the workbook to reflect the newly closed trade #89:
#  - Summary sheet aggregate metrics
#  - Strategy Status sheet row for MarketMaking
#  - All Trades sheet: append new trade row
#  - MarketMaking sheet: append new trade row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.69   # Current Capital
$wsSummary.Range("B4").Value = -0.32     # Total P&L $
$wsSummary.Range("B6").Value = 89        # Total Trades
$wsSummary.Range("B7").Value = 30        # Winning Trades
$wsSummary.Range("B9").Value = 33.71     # Win Rate %

# ---------------------------------------------------------------------
# Sheet 2: Strategy Status (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.69      # Capital
$wsStatus.Range("D4").Value = 89         # Trades
$wsStatus.Range("E4").Value = -0.32      # P&L $
$wsStatus.Range("F4").Value = -0.31      # P&L %
$wsStatus.Range("G4").Value = 33.71      # Win Rate %

# ---------------------------------------------------------------------
# Helper to append the new trade row (#89, spreadsheet row 90) to a
# trades-log style sheet (both "All Trades" and "MarketMaking" share
# the same layout / same new row).
# ---------------------------------------------------------------------
function Add-Trade89Row($ws) {
    $row = 90

    $ws.Cells.Item($row, 1).Value = 89

    # Date / Time must stay plain text, not get auto-converted to Excel
    # date/time serial values.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "15:53:49"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.65
    $ws.Cells.Item($row, 7).Value = 0.66
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 1.5385
    $ws.Cells.Item($row, 10).Value = 0.01
    $ws.Cells.Item($row, 11).Value = 99.69
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

# ---------------------------------------------------------------------
# Sheet 3: All Trades
# ---------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-Trade89Row $wsAllTrades

# ---------------------------------------------------------------------
# Sheet 4: MarketMaking
# ---------------------------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade89Row $wsMarketMaking
